# Atualizacao de bases das ligas (Croatia 3NL) - re-sync of several match
# rows whose underlying records were reordered/refreshed in the source feed.
# Each affected row below is overwritten, column by column (B:AD), with the
# full record that belongs there after the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($ws, $row, $vals) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $vals[$i]
    }
}

# Row 7 <= original row 9 data
$row7 = @(7126859, "Croatia 3NL", 45165.52083333334, "NK Granicar Zupanja", "Slavija Pleternica", 1, 1, 1, 1, "D", 1.5, 4, 5, 1.5, 4, 5, -1, 1.8, 2, 3, 1.85, 1.95, -1, 3, -1, -1, 1, -1, 0.95)
Set-RowValues $ws 7 $row7

# Row 8 <= original row 7 data
$row8 = @(7126857, "Croatia 3NL", 45165.52083333334, "NK Bistra", "NK Vrapce", 3, 2, 3, 0, "H", 2.875, 3.5, 2.1, 2.9, 3.5, 2.1, 0.25, 1.925, 1.875, 2.75, 1.825, 1.975, 1.9, -1, -1, 0.925, -1, 0.825, -1)
Set-RowValues $ws 8 $row8

# Row 9 <= original row 8 data
$row9 = @(7126858, "Croatia 3NL", 45165.52083333334, "Sava Strmec", "Lucko", 0, 2, 0, 1, "A", 2.75, 3.4, 2.2, 3.75, 3.6, 1.8, 0.5, 1.95, 1.85, 2.75, 1.9, 1.9, -1, -1, 0.8, -1, 0.8500000000000001, -1, 0.8999999999999999)
Set-RowValues $ws 9 $row9

# Row 33 <= original row 34 data
$row33 = @(7291473, "Croatia 3NL", 45203.47916666666, "NK Maksimir", "NK Mladost Petrinja", 5, 0, 4, 0, "H", 1.25, 6, 7, 1.25, 6, 7.5, -1.75, 1.9, 1.9, 3, 1.975, 1.825, 0.25, -1, -1, 0.8999999999999999, -1, 0.9750000000000001, -1)
Set-RowValues $ws 33 $row33

# Row 34 <= original row 33 data
$row34 = @(7291472, "Croatia 3NL", 45203.47916666666, "NK Lukavec", "Sava Strmec", 3, 0, 1, 0, "H", 2.2, 3.6, 2.6, 2.2, 3.6, 2.625, -0.25, 2, 1.8, 2.5, 1.8, 2, 1.2, -1, -1, 1, -1, 0.8, -1)
Set-RowValues $ws 34 $row34

# Row 78 <= original row 79 data
$row78 = @(7519478, "Croatia 3NL", 45256.40625, "NK Bistra", "Lucko", 1, 2, 0, 0, "A", 3, 3.6, 2, 3, 3.6, 2, 0.25, 2, 1.8, 2.5, 1.8, 2, -1, -1, 1, -1, 0.8, 0.8, -1)
Set-RowValues $ws 78 $row78

# Row 79 <= original row 78 data
$row79 = @(7519479, "Croatia 3NL", 45256.40625, "Sava Strmec", "NK Ponikve", 0, 1, 0, 1, "A", 2.75, 3.4, 2.2, 2.75, 3.4, 2.2, 0.25, 1.8, 2, 2.75, 1.9, 1.9, -1, -1, 1.2, -1, 1, -1, 0.8999999999999999)
Set-RowValues $ws 79 $row79

# Row 127 <= original row 129 data
$row127 = @(8163880, "Croatia 3NL", 45413.52083333334, "NK Dinamo Odranski Obre", "Sava Strmec", 2, 2, 1, 1, "D", 1.909, 3.4, 3.4, 1.8, 3.6, 3.75, -0.5, 1.825, 1.975, 3.25, 1.925, 1.875, -1, 2.6, -1, -1, 0.9750000000000001, 0.925, -1)
Set-RowValues $ws 127 $row127

# Row 129 <= original row 127 data
$row129 = @(8163883, "Croatia 3NL", 45413.52083333334, "NK Zelina", "NK Mladost Petrinja", 2, 1, 0, 1, "H", 2, 3.3, 3.25, 2.05, 3.5, 3, -0.25, 1.85, 1.95, 3.5, 1.825, 1.975, 1.05, -1, -1, 0.8500000000000001, -1, -1, 0.9750000000000001)
Set-RowValues $ws 129 $row129

# Row 136 <= original row 138 data
$row136 = @(8229446, "Croatia 3NL", 45430.52083333334, "RNK Split", "NK Omis", 2, 1, 1, 1, "H", 2.1, 3.4, 2.9, 1.8, 3.6, 3.6, -0.5, 1.85, 1.95, 2.75, 1.8, 2, 0.8, -1, -1, 0.8500000000000001, -1, 0.4, -0.5)
Set-RowValues $ws 136 $row136

# Row 137 <= original row 136 data
$row137 = @(8229445, "Croatia 3NL", 45430.52083333334, "NK Vodice", "NK Zadar", 1, 1, 0, 1, "D", 7, 4.5, 1.333, 7, 4.5, 1.333, 1.5, 1.9, 1.9, 3, 1.825, 1.975, -1, 3.5, -1, 0.8999999999999999, -1, -1, 0.9750000000000001)
Set-RowValues $ws 137 $row137

# Row 138 <= original row 137 data
$row138 = @(8229444, "Croatia 3NL", 45430.52083333334, "NK Neretva", "Zmaj Makarska", 1, 0, 1, 0, "H", 1.727, 3.75, 3.75, 1.5, 4.2, 4.75, -1, 1.8, 2, 3, 1.875, 1.925, 0.5, -1, -1, 0, 0, -1, 0.925)
Set-RowValues $ws 138 $row138
